$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the row above (row 5) down into the new row 6,
# matching the date/text/number styling already used for task rows.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B5:H5").Copy()
$ws.Range("B6:H6").PasteSpecial(-4122)

# Fill in the new weekly-report entry.
$ws.Range("A6").Value = 40918
$ws.Range("C6").Value = "1 phần SDD"
$ws.Range("D6").Value = "chưa hoàn thành"
$ws.Range("B6").Value = "viết các chức năng SDD như SRS"
$ws.Range("E6").Value = 4

# Match the workbook's saved selection state.
$ws.Range("H9").Select()
